# Add a "cfop" column to the "PI hours" sheet and a new "cfop hours" sheet
# summarising cfop totals, mirroring the existing "department hours" /
# "unit(accumulative) hours" layout.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "PI hours": append a new "cfop" column (G) after "app" (F)
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("PI hours")

# Copy the header cell's formatting (bold, centered, bordered) onto G1
# without introducing a brand-new style entry.
$ws1.Range("F1").Copy()
$ws1.Range("G1").PasteSpecial(-4122)
$ws1.Range("G1").Value = "cfop"

$ws1.Range("G2").Value = "['cfop_CHOUDHURY', 'cfop_RRC']"
$ws1.Range("G3").Value = "['cfop_NH']"
$ws1.Range("G4").Value = "['cfop_KWIAT']"
$ws1.Range("G5").Value = "['cfop_MITRA']"

# ------------------------------------------------------------------
# 2. New sheet "cfop hours" appended at the end of the workbook
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "cfop hours"

# Borrow header / index-column formatting from the existing
# "unit(accumulative) hours" sheet so the new sheet matches the others.
$ws3 = $wb.Worksheets.Item("unit(accumulative) hours")
$ws3.Range("B1:D1").Copy()
$ws4.Range("B1:D1").PasteSpecial(-4122)
$ws3.Range("A2").Copy()
$ws4.Range("A2:A6").PasteSpecial(-4122)

$ws4.Range("B1").Value = "cfop"
$ws4.Range("C1").Value = "hours"
$ws4.Range("D1").Value = "percentage"

$ws4.Range("A2").Value = 0
$ws4.Range("B2").Value = "cfop_RRC"
$ws4.Range("C2").Value = 100
$ws4.Range("D2").Value = 50.50505050505051

$ws4.Range("A3").Value = 1
$ws4.Range("B3").Value = "cfop_CHOUDHURY"
$ws4.Range("C3").Value = 50
$ws4.Range("D3").Value = 25.25252525252525

$ws4.Range("A4").Value = 2
$ws4.Range("B4").Value = "cfop_NH"
$ws4.Range("C4").Value = 35
$ws4.Range("D4").Value = 17.67676767676768

$ws4.Range("A5").Value = 3
$ws4.Range("B5").Value = "cfop_KWIAT"
$ws4.Range("C5").Value = 7
$ws4.Range("D5").Value = 3.535353535353535

$ws4.Range("A6").Value = 4
$ws4.Range("B6").Value = "cfop_MITRA"
$ws4.Range("C6").Value = 6
$ws4.Range("D6").Value = 3.03030303030303

# Keep "PI hours" as the active/selected tab, matching the original workbook.
$ws1.Activate()
